$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# C2 ("E200 DS_enhanced ..."): the run "(running)" has finished -> fill in
# the result and drop the green "still running" highlight color.
$cellC2 = $ws.Range("C2")
$oldRun = "J0: (running)"
$newRun = "J0: 27.8836 +/- 0.08"
$startPos = $cellC2.Text.IndexOf($oldRun) + 1
$runChars = $cellC2.Characters($startPos, $oldRun.Length)
$runChars.Text = $newRun
$newRunChars = $cellC2.Characters($startPos, $newRun.Length)
$newRunChars.Font.ColorIndex = -4105

# E2: simplify the "E200 NoA: " label to "E200: "
$cellE2 = $ws.Range("E2")
$cellE2.Replace("E200 NoA: ", "E200: ") | Out-Null

# Update the sheet's remembered selection
$ws.Range("A3").Select() | Out-Null

Write-Output "done"
